# edit.ps1
# Applies the weekly CompStat data refresh:
#  - Updates the "Volume / Number" run and the "Report Covering the Week" date
#    run inside their rich-text cells (A8, C9).
#  - Refreshes the Week-to-Date / 28-Day / Year-to-Date / %Chg figures in the
#    crime-category table (rows 14-30).
# Some cells flip between a numeric figure and the "0"/"***.*" placeholder
# text used for zero-complaint / undefined-percentage rows; for those we
# clone both the value AND the style from a same-shaped donor cell via
# Range.Copy (so the resulting style index / shared-string placeholder
# matches exactly), then overwrite the value when the target is numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Volume/Number run ("24" -> "25") ---
$ws.Range("A8").Value = "Volume 30   Number  25"

# --- Header text: week-covering date runs ---
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# --- Crime table: cells that change between numeric and placeholder-text styles ---
$ws.Range("C23").Copy($ws.Range("C15"))
$ws.Range("C36").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K36").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("C23").Copy($ws.Range("C22"))
$ws.Range("C23").Copy($ws.Range("C26"))
$ws.Range("C36").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("K36").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("C36").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("C36").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K36").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("C36").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K36").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("C23").Copy($ws.Range("G30"))
$ws.Range("E23").Copy($ws.Range("H30"))

# --- Crime table: straightforward numeric value refreshes ---
$ws.Range("N14").Value = -80
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 22.222222222222
$ws.Range("N15").Value = -26.666666666666
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = -11.538461538461
$ws.Range("I16").Value = 106
$ws.Range("J16").Value = 92
$ws.Range("K16").Value = 15.217391304347
$ws.Range("L16").Value = 107.843137254902
$ws.Range("M16").Value = -8.620689655172
$ws.Range("N16").Value = -76.906318082788
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -58.333333333333
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 138
$ws.Range("J17").Value = 141
$ws.Range("K17").Value = -2.127659574468
$ws.Range("L17").Value = 23.214285714285
$ws.Range("M17").Value = 15.966386554621
$ws.Range("N17").Value = 2.985074626865
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -71.428571428571
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 125
$ws.Range("K18").Value = -35.2
$ws.Range("L18").Value = -28.318584070796
$ws.Range("M18").Value = -59.5
$ws.Range("N18").Value = -91.588785046729
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 30.769230769230
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 24
$ws.Range("I19").Value = 319
$ws.Range("J19").Value = 299
$ws.Range("K19").Value = 6.688963210702
$ws.Range("L19").Value = 44.343891402714
$ws.Range("M19").Value = 64.432989690721
$ws.Range("N19").Value = 10.380622837370
$ws.Range("C20").Value = 10
$ws.Range("E20").Value = 42.857142857142
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 26.086956521739
$ws.Range("I20").Value = 166
$ws.Range("J20").Value = 122
$ws.Range("K20").Value = 36.065573770491
$ws.Range("L20").Value = 67.676767676767
$ws.Range("M20").Value = -4.597701149425
$ws.Range("N20").Value = -90.465249856404
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 150
$ws.Range("G21").Value = 160
$ws.Range("H21").Value = -6.25
$ws.Range("I21").Value = 823
$ws.Range("J21").Value = 790
$ws.Range("K21").Value = 4.177215189873
$ws.Range("L21").Value = 35.584843492586
$ws.Range("M21").Value = 1.479654747225
$ws.Range("N21").Value = -77.208529493215
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 125
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = -1.574803149606
$ws.Range("I24").Value = 611
$ws.Range("J24").Value = 670
$ws.Range("K24").Value = -8.805970149253
$ws.Range("L24").Value = -6.144393241167
$ws.Range("M24").Value = 16.826003824091
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = -29.166666666666
$ws.Range("I25").Value = 239
$ws.Range("J25").Value = 248
$ws.Range("K25").Value = -3.629032258064
$ws.Range("L25").Value = 15.458937198067
$ws.Range("M25").Value = -33.977900552486
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = 0
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 33
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 3.125
$ws.Range("L27").Value = 17.857142857142
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = -75
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -71.428571428571
